# Insert a new daily price record for "Macroferia Regional de Talca - Pepino
# ensalada" as row 541, shifting the existing rows 541:566 down to 542:567
# (dimension grows from A1:R566 to A1:R567).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 541 (and everything below it) down by one row.
$ws.Rows.Item(541).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A541").Value = 5
$ws.Range("B541").Value = "Macroferia Regional de Talca"
$ws.Range("C541").Value = "Maule"
$ws.Range("D541").Value = 44939
$ws.Range("E541").Value = 7
$ws.Range("F541").Value = 100112043
$ws.Range("G541").Value = "Pepino ensalada"
$ws.Range("H541").Value = "Sin especificar"
$ws.Range("I541").Value = "Primera"
$ws.Range("J541").Value = 500
$ws.Range("K541").Value = 8000
$ws.Range("L541").Value = 9000
$ws.Range("M541").Value = 8400
$ws.Range("N541").Value = "$/caja 80 unidades"
$ws.Range("O541").Value = "Región del Maule"
$ws.Range("P541").Value = 105
$ws.Range("Q541").Value = 80
$ws.Range("R541").Value = "Hortaliza"
